# feat: add 2022-Q1 data
#
# Before:  2021-Q1, 总计
# After:   2021-Q1, 2022-Q1 (new fund-holdings sheet), 总计 (updated summary)
#
# Strategy:
#  1. Duplicate the existing "总计" sheet (Copy) so the new "总计" keeps the
#     exact same cell formatting (header style, A-column index style, ...).
#  2. Rename the original "总计" sheet to "2022-Q1" and turn it into a
#     fund-holdings sheet (same shape as the "2021-Q1" sheet: columns
#     B..H = 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
#  3. Rename the duplicated sheet back to "总计" and insert a new first data
#     row for the "2022-Q1" summary line, pushing the existing "2021-Q1" row
#     down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Duplicate "总计" up front so formatting for the new summary sheet
#        comes from the real original (border/bold "总计" style), not from a
#        freshly minted blank sheet.
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item("总计 (2)")

# --- 2. Turn the original "总计" sheet into the "2022-Q1" holdings sheet.
$totalSheet.Name = "2022-Q1"
$ws = $totalSheet

# Make sure the header formatting (style "2" on B1) extends across the
# wider header row before filling it in.
$ws.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# Row 2 - the fund holding itself. A2 (index column) keeps its existing
# style; text columns that look numeric need an explicit text format so
# they stay strings instead of being coerced to numbers.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "000049"
$ws.Cells.Item(2, 3).Value = "中银标普全球精选自然资源等权重指数(QDII)"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "0.27"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "89.72"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "1.15"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "0.0031"
$ws.Cells.Item(2, 8).Value = 8

# --- 3. Turn the duplicated sheet into the new "总计" summary sheet.
$newTotalSheet.Name = "总计"
$ts = $newTotalSheet

# Insert a fresh row above the existing "2021-Q1" summary row, then clean
# up the blank cells the insert creates (they pick up stray formatting)
# before writing the new "2022-Q1" totals into them.
$ts.Rows.Item(2).Insert()
$ts.Range("B2:D2").ClearFormats()

$ts.Range("B1").Copy()
$ts.Range("A2").PasteSpecial(-4122)

$ts.Cells.Item(2, 1).Value = 0
$ts.Cells.Item(2, 2).Value = "2022-Q1"
$ts.Cells.Item(2, 3).Value = 1
$ts.Cells.Item(2, 4).Value = 0

# The old "2021-Q1" summary row got pushed down to row 3; its index needs
# to move from 0 to 1 to match its new position.
$ts.Cells.Item(3, 1).Value = 1

# Leave the originally-active sheet selected, same as before the edit.
$wb.Worksheets.Item("2021-Q1").Activate()
